$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
# Force a real run rewrite: the target text already equals the
# concatenation of the existing runs, so assigning it directly would be a
# no-op for the underlying run structure. Round-trip through a different
# value first so the merge into a single run actually takes effect.
$tr.Text = "__tmp__"
$tr.Text = "Below section-level"
